$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.671904563903809
$ws.Range("B1").Value = 3.963855504989624
$ws.Range("C1").Value = 7.640444755554199
$ws.Range("D1").Value = 7.806790351867676
$ws.Range("E1").Value = 5.962972640991211
